$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell B11 value from "R40" to "1".
# The target value "1" looks numeric, so to keep it stored as text (as in the
# source diff, which adds a new shared string "1" rather than a numeric value)
# we explicitly format the cell as Text before assigning the value. Otherwise
# Excel's automatic type inference would convert it to a number.
$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
